# Apply updated Betfair back/lay odds for 2026-01-01
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.74
$ws.Range("G2").Value = 1.76
$ws.Range("H2").Value = 4.6
$ws.Range("N2").Value = 7.2
$ws.Range("R2").Value = 1.86
$ws.Range("T2").Value = 1.5
$ws.Range("U2").Value = 2.86
$ws.Range("W2").Value = 2.3
$ws.Range("Y2").Value = 30
$ws.Range("Z2").Value = 44
$ws.Range("AB2").Value = 16.5
$ws.Range("AD2").Value = 18.5
$ws.Range("AE2").Value = 44
$ws.Range("AF2").Value = 15.5
$ws.Range("AJ2").Value = 20
$ws.Range("AN2").Value = 6

# Row 3
$ws.Range("F3").Value = 1.92
$ws.Range("G3").Value = 1.95
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 4.1
$ws.Range("N3").Value = 5.1
$ws.Range("R3").Value = 1.56
$ws.Range("S3").Value = 2.66
$ws.Range("W3").Value = 2.04
$ws.Range("AN3").Value = 9.4

# Row 4
$ws.Range("H4").Value = 4
$ws.Range("J4").Value = 3.95
$ws.Range("N4").Value = 4.7
$ws.Range("O4").Value = 1.21
$ws.Range("P4").Value = 2.32
$ws.Range("Q4").Value = 1.62
$ws.Range("R4").Value = 1.51
$ws.Range("S4").Value = 2.56
$ws.Range("U4").Value = 2.28
$ws.Range("Y4").Value = 990

# Row 5
$ws.Range("G5").Value = 2.88
$ws.Range("H5").Value = 2.98
$ws.Range("J5").Value = 2.74
$ws.Range("P5").Value = 1.54
$ws.Range("S5").Value = 3.1
$ws.Range("W5").Value = 1.53

# Row 6
$ws.Range("J6").Value = 3.55
$ws.Range("S6").Value = 2.66
$ws.Range("V6").Value = 1.04

# Row 7
$ws.Range("H7").Value = 1.25
$ws.Range("K7").Value = 8

# Row 8
$ws.Range("H8").Value = 1.09
$ws.Range("I8").Value = 4.8
$ws.Range("J8").Value = 2.88
$ws.Range("N8").Value = 1.63
$ws.Range("P8").Value = 1.63
$ws.Range("V8").Value = 1.26

# Row 9
$ws.Range("F9").Value = 2.26
$ws.Range("G9").Value = 2.28
$ws.Range("P9").Value = 1.84
$ws.Range("T9").Value = 1.88
$ws.Range("W9").Value = 1.78
$ws.Range("AH9").Value = 18.5

# Row 10
$ws.Range("F10").Value = 1.56
$ws.Range("G10").Value = 1.57
$ws.Range("P10").Value = 2.46
$ws.Range("Q10").Value = 1.66
$ws.Range("R10").Value = 1.57
$ws.Range("S10").Value = 2.7
$ws.Range("X10").Value = 23
$ws.Range("AA10").Value = 180
$ws.Range("AH10").Value = 20
$ws.Range("AL10").Value = 28
$ws.Range("AN10").Value = 6.8
$ws.Range("AO10").Value = 80

# Row 11
$ws.Range("G11").Value = 2.08
$ws.Range("H11").Value = 3.55
$ws.Range("S11").Value = 2.28
$ws.Range("W11").Value = 1.93
$ws.Range("Y11").Value = 25
$ws.Range("AB11").Value = 17
$ws.Range("AF11").Value = 18.5
$ws.Range("AJ11").Value = 28
$ws.Range("AK11").Value = 22
$ws.Range("AN11").Value = 10.5

# Row 12
$ws.Range("L12").Value = 1.41
$ws.Range("O12").Value = 1.32
$ws.Range("P12").Value = 2.02
$ws.Range("Q12").Value = 1.96
$ws.Range("S12").Value = 3.45
$ws.Range("X12").Value = 15.5
$ws.Range("AJ12").Value = 29
$ws.Range("AN12").Value = 17
$ws.Range("AO12").Value = 34

# Row 13
$ws.Range("F13").Value = 9
$ws.Range("G13").Value = 9.2
$ws.Range("M13").Value = 1.04
$ws.Range("O13").Value = 1.22
$ws.Range("P13").Value = 2.4
$ws.Range("Q13").Value = 1.7
$ws.Range("R13").Value = 1.55
$ws.Range("S13").Value = 2.74
$ws.Range("T13").Value = 1.96
$ws.Range("V13").Value = 3.35
$ws.Range("Z13").Value = 8.6
$ws.Range("AB13").Value = 32
$ws.Range("AH13").Value = 25
$ws.Range("AO13").Value = 5.7

Write-Host "Applied 103 cell updates"
